$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(-2.228, -0.443, -12.035, -0.12, 0.561, 0.68, -0.275, -0.102, 0.295, -0.027, -0.607, 2.643)
    3 = @(1.048, 10.552, -12.301, -1.157, 7.089, 9.049, 3.06, 17.618, 5.572, 7.425, 9.618, 7.468)
    4 = @(-10.816, 16.784, -46.319, 2.362, 1.931, 6.22, -15.483, 4.327, -13.922, 36.198, 8.641, 0.779)
    5 = @(-5.291, 24.672, -47.002, 10.657, 7.713, 13.16, -13.925, 14.21, -5.329, 58.04, 16.909, 6.44)
    6 = @(0.68, 143.816, -47.74, 150.799, 108.673, 42.893, -3.282, 62.688, 8.888, 222.65, 87.616, 7.892)
    7 = @(0.137, 19.689, -12.266, 20.372, 15.991, 7.463, -0.671, 10.312, 1.732, 26.645, 13.529, 1.544)
    8 = @(0.776, 1.008, 1.162, 0.949, 1.083, 1.024, 1.039, 0.813, 0.874, 1.015, 0.956, 0.567)
    9 = @(-6.226, 11.627, -21.448, 12.747, 7.386, -0.71, -8.954, 3.682, -5.344, 18.538, 5.849, -3.289)
    10 = @(-0.028, 0.9, -0.426, 0.807, 0.666, 0.28, -0.058, 0.446, 0.061, 0.944, 0.614, 0.086)
    11 = @(-0.007, 0.188, -0.111, 0.208, 0.141, 0.066, -0.013, 0.118, 0.012, 0.256, 0.134, 0.015)
    12 = @(-40.147, -30.196, -74.463, -38.017, -40.61, -36.068, -46.166, -35.435, -40.828, -29.459, -33.848, -28.439)
    13 = @(19.329, 21.132, 30.37, 24.406, 22.995, 24.221, 23.304, 21.603, 17.343, 27.527, 20.946, 10.09)
    14 = @(0.386, 0.708, 0.465, 0.46, 0.822, 0.517, 0.574, 0.386, 0.701, 0.626, 0.731, 0.328)
    15 = @(6.363, 8.062, 9.182, 7.626, 8.605, 8.173, 8.283, 6.63, 7.076, 8.107, 7.68, 4.832)
}

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M")

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "{0}{1}" -f $cols[$i], $r
        $ws.Range($addr).Value = $vals[$i]
    }
}
